# Relational Model hopefully works on
#
# 1. Rename the "Process ordering" use-case-name header to "makeSale".
# 2. Fix the "Costumer" -> "Customer" typo in the call-flow description.
# 3. Move the active selection to G2:H2 (the cell the author was working in).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "makeSale"
$ws.Range("G8").Value = "1.Customer makes a call and makes an order and specifies information (name,size,type) "

$ws.Activate()
$ws.Range("G2:H2").Select()
